# Applies the edits described by the diff to CryCompanywiseStockReport_1.xlsx
# Pattern summary:
#  1) A handful of row-pairs have their B (item code), C (item name), E (MRP),
#     F (qty) and G (value) contents swapped between the two rows (A and D
#     stay put).
#  2) A handful of individual rows have their quantity (F) reduced, with the
#     value (G) recomputed as Rate(D) * Qty(F).
#  3) The "Sub Total:" rows (column B) are recomputed as the sum of the G
#     column for the item rows belonging to that group.
#  4) The overall "Sub Total:"/"Grand Total:" rows (718/719) are recomputed
#     as the sum of every "Sub Total:" row in the sheet (excluding row 718
#     itself).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Cell($row1, $row2, $col) {
    $addr1 = "$col$row1"
    $addr2 = "$col$row2"
    $v1 = $ws.Range($addr1).Value2
    $v2 = $ws.Range($addr2).Value2
    $ws.Range($addr1).Value = $v2
    $ws.Range($addr2).Value = $v1
}

# ---------------------------------------------------------------------
# 1) Row-pair swaps: swap B, C, E, F, G between the two rows of each pair
# ---------------------------------------------------------------------
$pairs = @(
    @(161,162),
    @(290,291),
    @(292,293),
    @(297,298),
    @(304,305),
    @(308,309),
    @(317,318),
    @(479,480),
    @(506,507)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    foreach ($col in @('B','C','E','F','G')) {
        Swap-Cell $r1 $r2 $col
    }
}

# ---------------------------------------------------------------------
# 2) Simple quantity decreases: set new F (qty) and recompute G = D * F
# ---------------------------------------------------------------------
$qtyChanges = @{
    115 = 10
    184 = 55
    186 = 20
    217 = 48
    282 = 17
    310 = 14
    324 = 50
    370 = 230
    504 = 21
    531 = 219
    534 = 132
    555 = 6
    620 = 369
    625 = 325
    662 = 45
    674 = 873
}

foreach ($row in $qtyChanges.Keys) {
    $newQty = $qtyChanges[$row]
    $rate = $ws.Range("D$row").Value2
    $ws.Range("F$row").Value = $newQty
    $ws.Range("G$row").Value = [Math]::Round($rate * $newQty, 10)
}

# ---------------------------------------------------------------------
# 3) Recompute "Sub Total:" rows (B column) = SUM of G for the item rows
#    belonging to that group (rows between the previous text-label row
#    and this Sub Total row).
# ---------------------------------------------------------------------
$subtotalRows = @(123,193,218,295,328,372,508,535,556,628,668,680)

foreach ($subRow in $subtotalRows) {
    $r = $subRow - 1
    $sum = 0
    while ($r -ge 1) {
        $aVal = $ws.Range("A$r").Value2
        if ($aVal -ne $null -and $aVal.GetType().Name -eq "String") {
            # hit a text label (company header or a previous Sub Total row) -> stop
            break
        }
        $gVal = $ws.Range("G$r").Value2
        if ($gVal -ne $null -and $gVal -ne "") {
            $sum = $sum + $gVal
        }
        $r = $r - 1
    }
    $ws.Range("B$subRow").Value = [Math]::Round($sum, 10)
}

# ---------------------------------------------------------------------
# 4) Recompute the overall Sub Total (row 718) and Grand Total (row 719)
#    as the sum of every "Sub Total:" row on the sheet (excluding 718).
# ---------------------------------------------------------------------
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$grandSum = 0
for ($r = 1; $r -le $lastRow; $r++) {
    if ($r -eq 718) { continue }
    $label = $ws.Range("A$r").Value2
    if ($label -eq "Sub Total:") {
        $bVal = $ws.Range("B$r").Value2
        if ($bVal -ne $null -and $bVal -ne "") {
            $grandSum = $grandSum + $bVal
        }
    }
}

$grandSum = [Math]::Round($grandSum, 10)
$ws.Range("B718").Value = $grandSum
$ws.Range("B719").Value = $grandSum
